# sat base model til SGD ingen momentum
# Adds a new "Gamma" column (T) and four new result rows (14-17) to the
# CNN results log sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New column header: T1 = "Gamma" -------------------------------------
$ws.Range("T1").Value = "Gamma"
# Match the formatting of the other header cells (bold, centered, bordered)
$ws.Range("S1").Copy() | Out-Null
$ws.Range("T1").PasteSpecial(-4122) | Out-Null   # xlPasteFormats

# --- New data rows ---------------------------------------------------------
# Columns: A..T =
# Date&Time, Epochs, Batch size, Learning rate, Optimizer, Loss fn,
# Avg Time/Epoch, Image dim, Loss, Min Loss, Accuracy, Dataset, Device,
# Conv layers, Pools, Created by, Total training time, Weight decay,
# Learning rate decay, Gamma

$newRows = @(
    @("2024-1-4 21:19:16", 20, 64, [double]"1.000000000000001e-12", "SGD", "CEL", 35.9, 32, 1.7762, 1.2311, 45.3203, "FER2013", "cpu", 4, 2, "Stationær", 717.6847140835598, 0, $null, 0.1),
    @("2024-1-4 22:27:24", 20, 64, 0.001, "Adam", "CEL", 35.8, 32, 0.153, 0.07099999999999999, 97.405, "FER2013", "cpu", 4, 2, "Alfred", 716.4068887932226, 0, $null, $null),
    @("2024-1-5 10:31:45", 20, 64, 0.001, "Adam", "CEL", 37.2, 32, 0.214, 0.1574, 96.78149999999999, "FER2013", "cpu", 4, 2, "Alfred", 743.602515457198, 0, $null, $null),
    @("2024-1-5 10:45:27", 20, 64, 0.01, "SGD no momentum", "CEL", 35.9, 32, 1.3528, 1.3307, 51.4926, "FER2013", "cpu", 4, 2, "Alfred", 717.830478457734, 0.005, "None", $null)
)

$startRow = 14
for ($r = 0; $r -lt $newRows.Length; $r++) {
    $rowValues = $newRows[$r]
    $excelRow = $startRow + $r
    for ($c = 0; $c -lt $rowValues.Length; $c++) {
        $v = $rowValues[$c]
        if ($null -ne $v) {
            $ws.Cells.Item($excelRow, $c + 1).Value = $v
        }
    }
}
